$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 347
$ws.Range("F7").Value = 896
$ws.Range("F8").Value = 64
$ws.Range("F9").Value = 539
$ws.Range("F11").Value = 300
$ws.Range("F12").Value = 1169
$ws.Range("F15").Value = 42
$ws.Range("F17").Value = 6716
$ws.Range("F19").Value = 75
$ws.Range("F21").Value = 7635
$ws.Range("F24").Value = 3416
$ws.Range("F25").Value = 34
$ws.Range("F26").Value = 2145
$ws.Range("F27").Value = 916
$ws.Range("F28").Value = 4522
$ws.Range("F29").Value = 173
$ws.Range("F31").Value = 73
$ws.Range("F33").Value = 240
$ws.Range("F34").Value = 203
$ws.Range("F35").Value = 1767
$ws.Range("F37").Value = 195
$ws.Range("F39").Value = 5
$ws.Range("F40").Value = 19
$ws.Range("F41").Value = 1243
$ws.Range("F42").Value = 1861
$ws.Range("F43").Value = 2150
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 54
$ws.Range("F7").Value = 83
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1239
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1239
$ws.Range("F7").Value = 347
$ws.Range("F9").Value = 896
$ws.Range("F10").Value = 64
$ws.Range("F11").Value = 539
$ws.Range("F13").Value = 300
$ws.Range("F14").Value = 1169
$ws.Range("F18").Value = 42
$ws.Range("F20").Value = 6716
$ws.Range("F22").Value = 75
$ws.Range("F24").Value = 7635
$ws.Range("F27").Value = 3416
$ws.Range("F28").Value = 34
$ws.Range("F29").Value = 2145
$ws.Range("F30").Value = 916
$ws.Range("F31").Value = 4522
$ws.Range("F32").Value = 173
$ws.Range("F34").Value = 73
$ws.Range("F35").Value = 54
$ws.Range("F37").Value = 240
$ws.Range("F38").Value = 1767
$ws.Range("F40").Value = 195
$ws.Range("F42").Value = 5
$ws.Range("F43").Value = 19
$ws.Range("F44").Value = 1243
$ws.Range("F45").Value = 1861
$ws.Range("F47").Value = 2150
$ws.Range("F49").Value = 83
